# Add the new computer "esmith10laptop" as an extra column (H) in the
# ComputerFolders table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H1").Value = "esmith10laptop"
$ws.Range("H2").Value = "C:\E\Local\Lab\[07] Transcription\RawData"
$ws.Range("H3").Value = "C:\E\Local\Lab\[07] Transcription\FISHAnalysisData"
$ws.Range("H4").Value = "C:\E\SkyDrive\Lab\[07] Transcription\LivemRNAData"

# Resize the new columns (and previously-unsized column G) to fit their
# contents, similar to the bestFit column widths Excel recalculated for
# the rest of the table.
$ws.Columns("G").ColumnWidth = 38.1
$ws.Columns("H").ColumnWidth = 46.5

# Split the window so the new column is visible next to the row labels:
# column G becomes the first column of the right-hand pane.
$excel.ActiveWindow.SplitColumn = 6
$excel.ActiveWindow.SplitRow = 0

# Select the newly entered data for the new computer.
$ws.Range("H5").Select() | Out-Null
